# Update mapping new chips IDs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly scanned chip IDs for the mapping table (column A).
# The assignment order reproduces the order new shared-string table
# entries were created by the original author's edit.
$ws.Range("A15").Value = "103022550827"
$ws.Range("A17").Value = "987535843987"
$ws.Range("A16").Value = "786142748451"
$ws.Range("A18").Value = "796112435620"
$ws.Range("A19").Value = "1064020239013"
$ws.Range("A23").Value = "39911334407"
$ws.Range("A24").Value = "995451166718"
$ws.Range("A26").Value = "858500718555"
$ws.Range("A25").Value = "855499759477"
$ws.Range("A27").Value = "718631362413"

# Extend the print area to include the new rows and column C.
$ws.PageSetup.PrintArea = '$A$1:$C$30'

# Update the view: scroll so row 16 is at the top, and select A1:C30
# (mirrors the selection left by the author).
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:C30").Select()

$wb.Save()
